$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.572654128074646
$ws.Range("B1").Value = 3.672016143798828
$ws.Range("C1").Value = 5.51821231842041
$ws.Range("D1").Value = 1.370346784591675
$ws.Range("E1").Value = 0.7982211709022522
